# repull data, push all data, mean calculation
# Update the "dSF" (column F) values for several rows to reflect the
# repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -3
    6  = -6
    10 = -4
    12 = -2
    13 = 2
    14 = -1
    16 = -1
    18 = -2
    19 = 1
    23 = -4
    24 = 3
    25 = -1
    26 = 0
    27 = 1
    28 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
